$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4209.231
$ws.Range("J17").Value = 3417.6365
$ws.Range("L17").Value = 10252.9095
$ws.Range("N17").Value = -10588.9095
$ws.Range("H86").Value = 1841.1428
$ws.Range("J86").Value = 1944
$ws.Range("L86").Value = 1944
$ws.Range("N86").Value = -4190
$ws.Range("H89").Value = 1841.1428
$ws.Range("J89").Value = 1944
$ws.Range("L89").Value = 9720
$ws.Range("N89").Value = -20952
$ws.Range("H94").Value = 2660.75
$ws.Range("I94").Value = 2660.75
$ws.Range("K94").Value = 2660.75
$ws.Range("M94").Value = -2209.75
$ws.Range("H106").Value = 1820.1111
$ws.Range("I106").Value = 2300.4
$ws.Range("J106").Value = 1219.75
$ws.Range("K106").Value = 2300.4
$ws.Range("L106").Value = 1219.75
$ws.Range("M106").Value = -1669.4
$ws.Range("N106").Value = -2481.75
$ws.Range("H116").Value = 14551.556
$ws.Range("J116").Value = 3995.6667
$ws.Range("L116").Value = 3995.6667
$ws.Range("N116").Value = -10879.6667
$ws.Range("H138").Value = 2815.7964
$ws.Range("J138").Value = 2308.697
$ws.Range("L138").Value = 6926.091
$ws.Range("N138").Value = -17206.091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4434.0894
$ws.Range("I32").Value = 3593.077
$ws.Range("J32").Value = 15367.25
$ws.Range("K32").Value = 3593.077
$ws.Range("L32").Value = 15367.25
$ws.Range("M32").Value = -3306.077
$ws.Range("N32").Value = -15941.25
$ws.Range("H122").Value = 3145.7144
$ws.Range("I122").Value = 1210.5
$ws.Range("K122").Value = 3631.5
$ws.Range("M122").Value = -1181.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 79624.30499999999
$ws.Range("I86").Value = 3284.353
$ws.Range("K86").Value = 3284.353
$ws.Range("M86").Value = -2161.353
$ws.Range("H89").Value = 79624.30499999999
$ws.Range("I89").Value = 3284.353
$ws.Range("K89").Value = 16421.765
$ws.Range("M89").Value = -10805.765
$ws.Range("H94").Value = 833.0476
$ws.Range("I94").Value = 267.86667
$ws.Range("K94").Value = 267.86667
$ws.Range("M94").Value = 183.13333
$ws.Range("H99").Value = 1383
$ws.Range("I99").Value = 1283.4286
$ws.Range("J99").Value = 1499.1666
$ws.Range("K99").Value = 1283.4286
$ws.Range("L99").Value = 1499.1666
$ws.Range("M99").Value = 214.5714
$ws.Range("N99").Value = -4495.1666
$ws.Range("H134").Value = 5484.731
$ws.Range("J134").Value = 2659
$ws.Range("L134").Value = 7977
$ws.Range("N134").Value = -13047

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2374.3333
$ws.Range("I31").Value = 2277.3333
$ws.Range("K31").Value = 2277.3333
$ws.Range("M31").Value = -1982.3333
$ws.Range("H34").Value = 2374.3333
$ws.Range("I34").Value = 2277.3333
$ws.Range("K34").Value = 2277.3333
$ws.Range("M34").Value = -2075.3333
$ws.Range("H58").Value = 2901007.8
$ws.Range("I58").Value = 4833437
$ws.Range("K58").Value = 4833437
$ws.Range("M58").Value = -4833234
$ws.Range("H134").Value = 2565.1365
$ws.Range("I134").Value = 2226
$ws.Range("J134").Value = 5956.5
$ws.Range("K134").Value = 6678
$ws.Range("L134").Value = 17869.5
$ws.Range("M134").Value = -4143
$ws.Range("N134").Value = -22939.5
$ws.Range("H136").Value = 2901007.8
$ws.Range("I136").Value = 4833437
$ws.Range("K136").Value = 14500311
$ws.Range("M136").Value = -14497761

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10918.275
$ws.Range("I131").Value = 665.8
$ws.Range("J131").Value = 11719.25
$ws.Range("K131").Value = 1997.4
$ws.Range("L131").Value = 35157.75
$ws.Range("M131").Value = 3042.6
$ws.Range("N131").Value = -45237.75
$ws.Range("H139").Value = 5029.9
$ws.Range("I139").Value = 6417.476
$ws.Range("K139").Value = 19252.428
$ws.Range("M139").Value = -14112.428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 764.25
$ws.Range("J113").Value = 1224.7142
$ws.Range("L113").Value = 1224.7142
$ws.Range("N113").Value = -5564.7142
$ws.Range("H122").Value = 1089.3
$ws.Range("I122").Value = 889.2
$ws.Range("K122").Value = 2667.6
$ws.Range("M122").Value = -217.6000000000004
$ws.Range("H132").Value = 1375684.4
$ws.Range("I132").Value = 1924908.2
$ws.Range("J132").Value = 2624.75
$ws.Range("K132").Value = 5774724.6
$ws.Range("L132").Value = 7874.25
$ws.Range("M132").Value = -5772194.6
$ws.Range("N132").Value = -12934.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2785.4285
$ws.Range("I61").Value = 2571.7778
$ws.Range("J61").Value = 3170
$ws.Range("K61").Value = 2571.7778
$ws.Range("L61").Value = 3170
$ws.Range("M61").Value = -2369.7778
$ws.Range("N61").Value = -3574
$ws.Range("H93").Value = 18519076
$ws.Range("I93").Value = 559.2857
$ws.Range("J93").Value = 83333890
$ws.Range("K93").Value = 559.2857
$ws.Range("L93").Value = 83333890
$ws.Range("M93").Value = 688.7143
$ws.Range("N93").Value = -83336386
$ws.Range("H113").Value = 2785.4285
$ws.Range("I113").Value = 2571.7778
$ws.Range("J113").Value = 3170
$ws.Range("K113").Value = 2571.7778
$ws.Range("L113").Value = 3170
$ws.Range("M113").Value = -401.7777999999998
$ws.Range("N113").Value = -7510
$ws.Range("H132").Value = 2943.8
$ws.Range("I132").Value = 2015
$ws.Range("J132").Value = 3305
$ws.Range("K132").Value = 6045
$ws.Range("L132").Value = 9915
$ws.Range("M132").Value = -3515
$ws.Range("N132").Value = -14975

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1286.4
$ws.Range("I100").Value = 1108.25
$ws.Range("K100").Value = 2216.5
$ws.Range("M100").Value = -1675.5
$ws.Range("H122").Value = 22129.525
$ws.Range("I122").Value = 35362.39
$ws.Range("J122").Value = 1839.1333
$ws.Range("K122").Value = 106087.17
$ws.Range("L122").Value = 5517.3999
$ws.Range("M122").Value = -103637.17
$ws.Range("N122").Value = -10417.3999
$ws.Range("H132").Value = 2517.4211
$ws.Range("I132").Value = 1637.6428
$ws.Range("J132").Value = 4980.8
$ws.Range("K132").Value = 4912.928400000001
$ws.Range("L132").Value = 14942.4
$ws.Range("M132").Value = -2382.928400000001
$ws.Range("N132").Value = -20002.4
